$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Forecasted Consumption (MW)" values for rows 2..97 (column B)
$bValues = @(6250,6220,6210,6190,6160,6130,6120,6110,6100,6100,6100,6110,6130,6140,6160,6190,6250,6330,6440,6580,6750,6920,7080,7260,7460,7630,7780,7910,8030,8130,8190,8210,8210,8200,8180,8150,8100,8060,8030,7980,7930,7900,7890,7880,7850,7850,7850,7850,7910,7930,7930,7910,7880,7850,7850,7850,7890,7910,7930,7960,8000,8040,8080,8130,8170,8230,8290,8360,8450,8500,8510,8520,8500,8480,8440,8400,8340,8290,8230,8160,8040,7920,7820,7690,7520,7370,7200,7060,6940,6810,6690,6580,6500,6440,6390,6330)

$startRow = 2
$endRow = 97

for ($r = $startRow; $r -le $endRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)

    # Shift the timestamp forward by 2 days, keeping the time-of-day fraction
    $aCell.Value2 = $aCell.Value2 + 2

    # Update the forecasted consumption value
    $bCell.Value2 = $bValues[$r - $startRow]
}

Write-Host "Updated rows $startRow to $endRow"
